{"js": "// Applies the RQMT_game.docx edit (\"solved the lightning while animating problem\"):\n//  1. Add a trailing space run after \" { progression, pause, time display, time up }\"\n//     in the \"Time limit\" bullet.\n//  2. Replace the red \"DESTRUCTION -- ... { Big lightning strikes when sending }\" bullet\n//     with a rewritten, green version that adds proper grammar markers around\n//     \"{ Big\" / \"lightning strikes when sending }\" and appends a bold \"DONE!\".\n//  3. Insert a new \"Pause game\" bullet right after \"HELP page, menus\".\n\nasync function findParagraphByText(context, text) {\n  const body = context.document.body;\n  const paragraphs = body.paragraphs;\n  paragraphs.load('items/text');\n  await context.sync();\n  for (const p of paragraphs.items) {\n    if (p.text === text) {\n      return p;\n    }\n  }\n  return null;\n}\n\n// --- 1) \"Time limit\" bullet: append a trailing space run. ---\nconst timeLimitPara = await findParagraphByText(\n  context,\n  'Time limit { progression, pause, time display, time up }'\n);\nif (!timeLimitPara) {\n  throw new Error('Could not locate the \"Time limit\" paragraph.');\n}\ntimeLimitPara.insertText(' ', Word.InsertLocation.end);\nawait context.sync();\n\n// --- 2) DESTRUCTION bullet: insert the rewritten green paragraph before the\n//        old red one (as its own new paragraph), then delete the old one. ---\nconst destructionPara = await findParagraphByText(\n  context,\n  'DESTRUCTION -- when tiles disappear (bomb, arrow, send), some lightning appears in place for a bit { Big lightning strikes when sending }'\n);\nif (!destructionPara) {\n  throw new Error('Could not locate the \"DESTRUCTION\" paragraph.');\n}\n\nconst newDestructionPara = destructionPara.insertParagraph('', Word.InsertLocation.before);\nawait context.sync();\n\nconst destructionOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"5\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>DESTRUCTION -- when tiles disappear (bomb, arrow, send), some lightning appears in place for a bit</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>{ Big</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> lightning strikes when sending }</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>DONE!</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nnewDestructionPara.getRange('Whole').insertOoxml(destructionOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\ndestructionPara.delete();\nawait context.sync();\n\n// --- 3) Insert a new \"Pause game\" bullet right after \"HELP page, menus\". ---\nconst helpPagePara = await findParagraphByText(context, 'HELP page, menus');\nif (!helpPagePara) {\n  throw new Error('Could not locate the \"HELP page, menus\" paragraph.');\n}\nhelpPagePara.insertParagraph('Pause game', Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Applies the RQMT_game.docx edit (\"solved the lightning while animating problem\"):\n#  1. Add a trailing space run after \" { progression, pause, time display, time up }\"\n#     in the \"Time limit\" bullet.\n#  2. Replace the red \"DESTRUCTION -- ... { Big lightning strikes when sending }\" bullet\n#     with a rewritten, green version that adds proper grammar markers around\n#     \"{ Big\" / \"lightning strikes when sending }\" and appends a bold \"DONE!\".\n#  3. Insert a new \"Pause game\" bullet right after \"HELP page, menus\".\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, $text) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -eq ($text + \"`r\")) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- 1) \"Time limit\" bullet: append a trailing space run. ---\n$timeLimit = Find-ParagraphByText $d \"Time limit { progression, pause, time display, time up }\"\nif ($timeLimit -eq $null) {\n    throw \"Could not locate the 'Time limit' paragraph.\"\n}\n$endRange = $d.Range($timeLimit.Range.End - 1, $timeLimit.Range.End - 1)\n$endRange.InsertAfter(\" \") | Out-Null\n\n# --- 2) DESTRUCTION bullet: insert the rewritten green paragraph before the\n#        old red one (as its own new paragraph), then delete the old one. ---\n$destruction = Find-ParagraphByText $d \"DESTRUCTION -- when tiles disappear (bomb, arrow, send), some lightning appears in place for a bit { Big lightning strikes when sending }\"\nif ($destruction -eq $null) {\n    throw \"Could not locate the 'DESTRUCTION' paragraph.\"\n}\n\n$destruction.Range.InsertParagraphBefore() | Out-Null\n# $destruction now refers to the freshly-inserted EMPTY paragraph; the old\n# (red) paragraph got pushed to be its Next() sibling.\n$oldDestruction = $destruction.Next()\n\n$destructionXml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"5\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>DESTRUCTION -- when tiles disappear (bomb, arrow, send), some lightning appears in place for a bit</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>{ Big</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> lightning strikes when sending }</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:color w:val=\"00B050\"/>\n              </w:rPr>\n              <w:t>DONE!</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n$destruction.Range.InsertXML($destructionXml) | Out-Null\n\n$oldDestruction.Range.Delete() | Out-Null\n\n# --- 3) Insert a new \"Pause game\" bullet right after \"HELP page, menus\". ---\n$helpPage = Find-ParagraphByText $d \"HELP page, menus\"\nif ($helpPage -eq $null) {\n    throw \"Could not locate the 'HELP page, menus' paragraph.\"\n}\n$helpPage.Range.InsertParagraphAfter() | Out-Null\n$pauseGame = $helpPage.Next()\n$pauseGame.Range.InsertAfter(\"Pause game\") | Out-Null\n"}
